$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins / Losses / Ties) in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, thin box border,
# centered horizontally, aligned to top vertically)
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font().Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders().LineStyle = 1       # xlContinuous

# Fill in the team record (Wins=68, Losses=94, Ties=0) for every data row (2-45)
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = 68   # column AD
    $ws.Cells.Item($row, 31).Value = 94   # column AE
    $ws.Cells.Item($row, 32).Value = 0    # column AF
}
